# Apply the "remove buster support, preparing for ma5" edit described by the
# supplied diff to the downloads.xlsx workbook (Sheet 1).
#
# Summary of the substantive (data) changes:
#   - Row 37 height: 30 -> 15
#   - Row 115 gains an E value (114) which ripples B115/C115
#   - Rows 116-133 (months 2023/01 .. 2024/06) gain B/C/D/E data+formulas
#     (G values already existed in the original file)
#   - Selection changes to B116:B133 (active cell B116)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 37 formatting: shrink back down to the standard row height.
# ---------------------------------------------------------------------
$ws.Rows.Item(37).RowHeight = 15

# ---------------------------------------------------------------------
# 2. Row 115: add the missing "download count (month)" figure (column E).
#    B115/C115 already carry formulas, so they recompute automatically.
# ---------------------------------------------------------------------
$ws.Cells.Item(115, 5).Value = 114

# ---------------------------------------------------------------------
# 3. Rows 116-133: fill in B (total), C (download total), D (usb total),
#    E (download/month) following the same formula pattern used by every
#    preceding row. G (docker cumulative) is already populated.
# ---------------------------------------------------------------------
$monthData = @{
    116 = 188
    117 = 114
    118 = 263
    119 = 216
    120 = 263
    121 = 267
    122 = 234
    123 = 112
    124 = 166
    125 = 295
    126 = 176
    127 = 108
    128 = 105
    129 = 125
    130 = 327
    131 = 203
    132 = 235
}

# G (docker cumulative) was not yet populated for 2024/05-2024/06 (rows
# 131-132) nor 2024/06... actually rows 131-133 (2024/04..2024/06).
$dockerData = @{
    131 = 777
    132 = 823
    133 = 834
}

for ($r = 116; $r -le 133; $r++) {
    $ws.Cells.Item($r, 2).Formula = "=C$r+D$r+G$r"
    $ws.Cells.Item($r, 3).Formula = "=C" + ($r - 1) + "+E$r"
    $ws.Cells.Item($r, 4).Formula = "=D" + ($r - 1) + "+F$r"
    if ($monthData.ContainsKey($r)) {
        $ws.Cells.Item($r, 5).Value = $monthData[$r]
    }
    if ($dockerData.ContainsKey($r)) {
        $ws.Cells.Item($r, 7).Value = $dockerData[$r]
    }
}

# ---------------------------------------------------------------------
# 4. Update the sheet selection to match the edited range.
# ---------------------------------------------------------------------
$ws.Range("B116:B133").Select()

# ---------------------------------------------------------------------
# 5. Cosmetic: localize the built-in "Normal" cell style name to the
#    Japanese equivalent used in the edited workbook. (Best effort --
#    the hosted Style object only exposes a read-only Name in this
#    runtime, so this may be a no-op.)
# ---------------------------------------------------------------------
try {
    $wb.Styles.Item(1).Name = "標準"
} catch {
}
